$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Block 1: rows 34-44 ----
$block1 = @(
    @(34, 'Checking Ingress availability', 56),
    @(35, 'Adding case ibm-apiconnect', 7),
    @(36, 'Adding case ibm-mq', 4),
    @(37, 'Adding case ibm-licensing', 3),
    @(38, 'Adding case ibm-cp-common-services', 9),
    @(39, 'Creation of openshift-cert-manager-operator operator', 55),
    @(40, 'Creation of ibm-licensing-operator-app operator', 51),
    @(41, 'Creation of ibm-common-service-operator operator', 41),
    @(42, 'Creation of ibm-apiconnect operator', 51),
    @(43, 'Creation of ibm-mq operator', 55),
    @(44, 'Creation of commonservice instance', 1)
)
foreach ($item in $block1) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
}
$ws.Range("E34:E44").Formula = "=B34/86400"
$ws.Range("E34:E44").NumberFormat = "[h]:mm:ss;@"

# ---- Row 45 (standalone formula, not part of the shared group) ----
$ws.Cells.Item(45, 1).Value = 'Creation of APIConnectCluster instance'
$ws.Cells.Item(45, 2).Value = 1380
$ws.Range("E45").Formula = "=B45/86400"
$ws.Range("E45").NumberFormat = "[h]:mm:ss;@"

# ---- Row 46 (block 1 total) ----
$ws.Range("B46").Formula = "=SUM(B34:B45)"
$ws.Range("B46:D46").Font.Bold = $true
$ws.Range("E46").Formula = "=B46/86400"
$ws.Range("E46").Font.Bold = $true
$ws.Range("E46").NumberFormat = "[h]:mm:ss;@"

# ---- Row 32 label (added after block 1 was filled in, per shared-string order) ----
$ws.Cells.Item(32, 1).Value = '2023.4 2 january 2024'

# ---- Row 50: date number ----
$ws.Cells.Item(50, 1).Value = 20240116

# ---- Block 2: rows 51-82 ----
$block2 = @(
    @(51, 'Checking Ingress availability', 1005),
    @(52, 'Adding case ibm-integration-platform-navigator', 4),
    @(53, 'Adding case ibm-appconnect', 1),
    @(54, 'Adding case ibm-apiconnect', 3),
    @(55, 'Adding case ibm-cp-common-services', 3),
    @(56, 'Adding case ibm-eventendpointmanagement', 1),
    @(57, 'Adding case ibm-eventprocessing', 2),
    @(58, 'Adding case ibm-eventstreams', 1),
    @(59, 'Adding case ibm-eventautomation-flink', 2),
    @(60, 'Adding case ibm-licensing', 2),
    @(61, 'Adding case ibm-mq', 2),
    @(62, 'Creation of openshift-cert-manager-operator operator', 98),
    @(63, 'Creation of ibm-licensing-operator-app operator', 71),
    @(64, 'Creation of ibm-common-service-operator operator', 47),
    @(65, 'Creation of datapower-operator operator', 47),
    @(66, 'Creation of ibm-integration-platform-navigator operator', 49),
    @(67, 'Creation of ibm-appconnect operator', 46),
    @(68, 'Creation of ibm-apiconnect operator', 52),
    @(69, 'Creation of ibm-eventendpointmanagement operator', 95),
    @(70, 'Creation of ibm-eventautomation-flink.v1.1.1 operator', 65),
    @(71, 'Creation of ibm-eventprocessing.v1.1.1 operator', 68),
    @(72, 'Creation of ibm-eventstreams operator', 188),
    @(73, 'Creation of ibm-mq operator', 56),
    @(74, 'Creation of Dashboard instance', 119),
    @(75, 'Creation of DesignerAuthoring instance', 5925),
    @(76, 'Creation of APIConnectCluster instance', 1407),
    @(77, 'Creation of EventStreams instance', 652),
    @(78, 'Creation of EventEndpointManagement instance', 9),
    @(79, 'Creation of EventGateway instance', 7),
    @(80, 'Creation of PersistentVolumeClaim instance', 70),
    @(81, 'Creation of FlinkDeployment instance', 139),
    @(82, 'Creation of EventProcessing instance', 264)
)
foreach ($item in $block2) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
}

# ---- Row 83 (block 2 total, formula SUM) ----
$ws.Range("B83").Formula = "=SUM(B51:B82)"

# Apply the shared division formula across the WHOLE E51:E83 range at once
# (matches how row 83 ends up sharing si=3 with the rest, per the original file)
$ws.Range("E51:E83").Formula = "=B51/86400"
$ws.Range("E51:E83").NumberFormat = "[h]:mm:ss;@"

# ---- Highlight (red bold) E51, E75, E76 ----
foreach ($cellref in @("E51","E75","E76")) {
    $ws.Range($cellref).Font.Bold = $true
    $ws.Range($cellref).Font.Color = 255
}

# ---- Row 83 totals row bold styling ----
$ws.Range("B83:D83").Font.Bold = $true
$ws.Range("E83").Font.Bold = $true

# ---- Selection / view ----
$ws.Range("E77").Select()
